# Adding the changes we made on may 9th
# Appends 10 new accelerometer sample rows (x, y, z) to Sheet1, extending
# the used range from A1:C21 to A1:C31.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.443913996219635, 6.598052978515625, 6.372174263000488),
    @(5.017239570617676, -37.66963577270508, 0.8286991119384766),
    @(-11.97968292236328, -29.96435546875, -18.5664291381836),
    @(1.253250598907471, 2.138059854507446, 19.65131759643555),
    @(-0.3477138876914978, -27.24154853820801, 47.39388656616211),
    @(-66.68949127197266, -16.29855728149414, -26.75592613220215),
    @(-5.153336048126221, 1.702579855918884, 10.3774242401123),
    @(-7.020784378051758, -33.00114822387695, 7.655538558959961),
    @(6.488995552062988, -21.21427917480469, 29.90091133117676),
    @(0.7014479041099548, 5.056992530822754, 16.41212844848633)
)

$startRow = 22
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
